$d = $word.ActiveDocument

# Locate the unique "Developed a JavaScript" phrase (Algorithmic Trading Bot bullet)
# and bold only the "JavaScript" part of it.
$rng = $d.Content
$found = $rng.Find.Execute("Developed a JavaScript", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $jsStart = $rng.End - 10
    $jsRange = $d.Range($jsStart, $rng.End)
    $jsRange.Font.Bold = 1
    $jsRange.Font.BoldBi = 1
}
